$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.303.90"
$ws.Range("D3").Value = "3.494.62"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'588.63"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'133.96"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("E9").Value = "  +5.84%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "'0.390"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("D12").Value = "4.089.12"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "3.496.33"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "64.317.80"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'25.42"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "'388.44"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D23").Value = "3.633.95"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'74.20"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.51"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.26"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  +4.33%  "
$ws.Range("D34").Value = "3.523.17"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D36").Value = "'23.34"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'5.33"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "'165.64"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "'24.50"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").Value = "2.397.89"
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("E51").Value = "  -0.36%  "
